$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptos data from the scraped source (GitHub Actions refresh).
# Column D values are numeric-looking text (e.g. thousand-separated prices, tiny
# decimals) that must stay as literal text, not be reinterpreted as numbers, so
# we force the Text number format before assigning, then restore the default style.

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '69.388.93'
$c.Style = "Normal"
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '3.769.27'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -0.52%  '
$ws.Range('E4').Value = '  -0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '614.28'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -1.01%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '178.26'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.63%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '3.766.98'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -0.63%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('E9').Value = '  -1.29%  '
$ws.Range('E10').Value = '  -2.47%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '6.60'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +4.95%  '
$ws.Range('E12').Value = '  -1.45%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '39.92'
$c.Style = "Normal"
$ws.Range('E13').Value = '  -1.71%  '
$ws.Range('E14').Value = '  -3.38%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '4.395.22'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -0.86%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '3.767.98'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -1.02%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '69.459.54'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -0.76%  '
$ws.Range('E18').Value = '  -0.37%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.120'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -3.30%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '507.89'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -0.24%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '16.35'
$c.Style = "Normal"
$ws.Range('E21').Value = '  -3.09%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '9.41'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -1.37%  '
$ws.Range('E23').Value = '  +0.33%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.52'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +0.32%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '86.53'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -1.35%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '12.84'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -2.93%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '0.0000135'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -5.06%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '10.60'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -4.03%  '
$ws.Range('E29').Value = '  +0.05%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.58%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '2.98'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +3.52%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '30.67'
$c.Style = "Normal"
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('E34').Value = '  -0.68%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('E36').Value = '  -2.32%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '6.13'
$c.Style = "Normal"
$ws.Range('E37').Value = '  -1.36%  '
$ws.Range('E38').Value = '  +3.78%  '
$ws.Range('E39').Value = '  +2.70%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '457.35'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +9.15%  '
$ws.Range('E41').Value = '  -2.43%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '49.79'
$c.Style = "Normal"
$ws.Range('E42').Value = '  -2.51%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +4.90%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '44.57'
$c.Style = "Normal"
$ws.Range('E44').Value = '  -2.14%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '8.60'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -1.62%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.959.95'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -2.66%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '27.32'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -0.14%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '138.99'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +0.70%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '2.48'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +0.49%  '
